$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = [double]"9.670175771248068e-11"
$ws.Range("H2").Value = [double]"2.536439546556871e-10"
$ws.Range("I2").Value = [double]"0.1614735695256511"
$ws.Range("K2").Value = [double]"40.83887586882928"
$ws.Range("L2").Value = "[28.019713082071235, 53.65803865558732]"
$ws.Range("M2").Value = [double]"3.596294817498347e-09"
$ws.Range("N2").Value = [double]"3.596294817498347e-09"
$ws.Range("O2").Value = [double]"1.427710775505271"
$ws.Range("P2").Value = "[1.0629212381515014, 1.7925003128590404]"
$ws.Range("Q2").Value = [double]"1.768585278227874e-12"
$ws.Range("R2").Value = [double]"3.537170556455749e-12"
$ws.Range("S2").Value = [double]"55.71614755995852"
$ws.Range("T2").Value = "[48.018004424007415, 63.41429069590963]"
$ws.Range("W2").Value = [double]"17.001001001001"
$ws.Range("X2").Value = [double]"15.72372372372372"
$ws.Range("Y2").Value = [double]"18.27827827827828"

# Row 3 updates
$ws.Range("E3").Value = [double]"22.55000000000009"
$ws.Range("G3").Value = [double]"1.887379141862766e-14"
$ws.Range("H3").Value = [double]"2.34000852882533e-13"
$ws.Range("K3").Value = [double]"48.38812374603071"
$ws.Range("L3").Value = "[34.25617255717956, 62.52007493488187]"
$ws.Range("M3").Value = [double]"3.216791277793618e-10"
$ws.Range("N3").Value = [double]"6.433582555587236e-10"
$ws.Range("O3").Value = [double]"1.201289683354656"
$ws.Range("P3").Value = "[0.8868159442565782, 1.515763422452733]"
$ws.Range("Q3").Value = [double]"4.89697171701664e-12"
$ws.Range("R3").Value = [double]"4.89697171701664e-12"
$ws.Range("S3").Value = [double]"66.58854896626283"
$ws.Range("T3").Value = "[58.86141453182529, 74.31568340070038]"
$ws.Range("W3").Value = [double]"18.23863863863871"
$ws.Range("X3").Value = [double]"17.11001001001008"
$ws.Range("Y3").Value = [double]"19.36726726726734"
